# Fixed error in dictionary provided in the TestCases.xlsx file.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("StrategyDictionaries")

# ---------------------------------------------------------------------------
# 1) Fix the UtimateScalper dictionary string on StrategyDictionaries!B4:
#    python-style single quotes -> valid JSON double quotes.
# ---------------------------------------------------------------------------
$ws3.Range("B4").Value = '{"EMA_Fast": 9, "EMA_Slow": 55, "EMA_Trend": 200, "RSI": 2, "RSI_Low": 48, "RSI_High": 52, "ADX": 14, "ADX_Threshold": 0, "MACD_Fast": 12, "MACD_Slow": 26, "MACD_Signal": 9, "BB_Length": 34, "BB_Mult": 1}'

# ---------------------------------------------------------------------------
# 2) Delete the duplicate second test-case row (row 3) on Sheet1. This shifts
#    every row below it up by one (old row 10 -> 9, ..., old row 21 removed).
# ---------------------------------------------------------------------------
$ws1.Rows.Item(3).Delete()

# Deleting a row that sits above "whole column" data-validation ranges makes
# the engine shrink the bottom boundary of those ranges by one row (a quirk
# of the row-shift implementation). Restore them to span down to the very
# last row again, matching the original file exactly.
$vA = $ws1.Range("A2:A1048576").Validation
$vA.Delete()
$vA.Add(1, 1, 1, "1", "99999")
$vA.IgnoreBlank = $true
$vA.ShowInput = $true
$vA.ShowError = $true

$vGH = $ws1.Range("G2:H1048576").Validation
$vGH.Delete()
$vGH.Add(2, 1, 1, "0", "1000")
$vGH.IgnoreBlank = $true
$vGH.ShowInput = $true
$vGH.ShowError = $true

$vDE = $ws1.Range("D2:E1048576").Validation
$vDE.Delete()
$vDE.Add(4, 1, 1, "36526", "47484")
$vDE.IgnoreBlank = $true
$vDE.ShowInput = $true
$vDE.ShowError = $true

$vK = $ws1.Range("K2:K1048576").Validation
$vK.Delete()
$vK.Add(0, 1, 1, "", "")
$vK.InputTitle = "Optional Strategy Settings"
$vK.InputMessage = "Format  { ""key1"": value1, ""key2"": value2, … }"
$vK.IgnoreBlank = $true
$vK.ShowInput = $true
$vK.ShowError = $true

# ---------------------------------------------------------------------------
# 3) Update sheet selections/scroll position to match the saved view state.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("K2").Select()
$excel.ActiveWindow.ScrollColumn = 2

$ws3.Activate()
$ws3.Range("B12").Select()

$ws1.Activate()

# ---------------------------------------------------------------------------
# 4) Resize the workbook window.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Width = 29790
$excel.ActiveWindow.Height = 11685
